# Auto-generated cell value updates applying the Anima_Profits.xlsx diff
# (scheduled-runner refresh of scraped FFXIV market data; plain numeric values, no formulas)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 805.55554
$ws.Range("I6").Value = 100
$ws.Range("K6").Value = 300
$ws.Range("M6").Value = -188
$ws.Range("H8").Value = 271
$ws.Range("I8").Value = 209
$ws.Range("J8").Value = 550
$ws.Range("K8").Value = 627
$ws.Range("L8").Value = 1650
$ws.Range("M8").Value = -488
$ws.Range("N8").Value = -1928
$ws.Range("H116").Value = 3480
$ws.Range("I116").Value = 2400
$ws.Range("J116").Value = 6000
$ws.Range("K116").Value = 2400
$ws.Range("L116").Value = 6000
$ws.Range("M116").Value = 1042
$ws.Range("N116").Value = -12884
$ws.Range("H132").Value = 5646.722
$ws.Range("I132").Value = 5390.706
$ws.Range("J132").Value = 9999
$ws.Range("K132").Value = 16172.118
$ws.Range("L132").Value = 29997
$ws.Range("M132").Value = -13642.118
$ws.Range("N132").Value = -35057
$ws.Range("H137").Value = 2493.3667
$ws.Range("I137").Value = 2180.3
$ws.Range("J137").Value = 3119.5
$ws.Range("K137").Value = 6540.900000000001
$ws.Range("L137").Value = 9358.5
$ws.Range("M137").Value = -3990.900000000001
$ws.Range("N137").Value = -14458.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2351.3948
$ws.Range("I61").Value = 2013.7097
$ws.Range("J61").Value = 3846.8572
$ws.Range("K61").Value = 2013.7097
$ws.Range("L61").Value = 3846.8572
$ws.Range("M61").Value = -1801.7097
$ws.Range("N61").Value = -4270.8572
$ws.Range("H74").Value = 2610.476
$ws.Range("I74").Value = 1993.3572
$ws.Range("J74").Value = 3844.7144
$ws.Range("K74").Value = 1993.3572
$ws.Range("L74").Value = 3844.7144
$ws.Range("M74").Value = -1119.3572
$ws.Range("N74").Value = -5592.7144
$ws.Range("H77").Value = 2610.476
$ws.Range("I77").Value = 1993.3572
$ws.Range("J77").Value = 3844.7144
$ws.Range("K77").Value = 9966.786
$ws.Range("L77").Value = 19223.572
$ws.Range("M77").Value = -5598.786
$ws.Range("N77").Value = -27959.572
$ws.Range("H122").Value = 252559
$ws.Range("J122").Value = 3412
$ws.Range("L122").Value = 10236
$ws.Range("N122").Value = -15136
$ws.Range("H136").Value = 2351.3948
$ws.Range("I136").Value = 2013.7097
$ws.Range("J136").Value = 3846.8572
$ws.Range("K136").Value = 6041.1291
$ws.Range("L136").Value = 11540.5716
$ws.Range("M136").Value = -3491.1291
$ws.Range("N136").Value = -16640.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 28573838
$ws.Range("I20").Value = 2238.7144
$ws.Range("J20").Value = 71431240
$ws.Range("K20").Value = 2238.7144
$ws.Range("L20").Value = 71431240
$ws.Range("M20").Value = -1991.7144
$ws.Range("N20").Value = -71431734
$ws.Range("H134").Value = 2747.25
$ws.Range("I134").Value = 1814.0454
$ws.Range("J134").Value = 6169
$ws.Range("K134").Value = 5442.1362
$ws.Range("L134").Value = 18507
$ws.Range("M134").Value = -2907.1362
$ws.Range("N134").Value = -23577

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 7199.4
$ws.Range("I134").Value = 7436.5
$ws.Range("J134").Value = 6251
$ws.Range("K134").Value = 22309.5
$ws.Range("L134").Value = 18753
$ws.Range("M134").Value = -19774.5
$ws.Range("N134").Value = -23823

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 1823.75
$ws.Range("I31").Value = 3000
$ws.Range("K31").Value = 9000
$ws.Range("M31").Value = -8712
$ws.Range("H131").Value = 1108.6285
$ws.Range("J131").Value = 1157.1936
$ws.Range("L131").Value = 3471.5808
$ws.Range("N131").Value = -13551.5808

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5055.0156
$ws.Range("I70").Value = 5023.3657
$ws.Range("J70").Value = 5109.0835
$ws.Range("K70").Value = 5023.3657
$ws.Range("L70").Value = 5109.0835
$ws.Range("M70").Value = -4753.3657
$ws.Range("N70").Value = -5649.0835
$ws.Range("H73").Value = 5055.0156
$ws.Range("I73").Value = 5023.3657
$ws.Range("J73").Value = 5109.0835
$ws.Range("K73").Value = 5023.3657
$ws.Range("L73").Value = 5109.0835
$ws.Range("M73").Value = -4087.3657
$ws.Range("N73").Value = -6981.0835
$ws.Range("H113").Value = 60604.293
$ws.Range("I113").Value = 72982.86
$ws.Range("K113").Value = 72982.86
$ws.Range("M113").Value = -70812.86
$ws.Range("H132").Value = 3241.7778
$ws.Range("I132").Value = 2855.138
$ws.Range("J132").Value = 4843.5713
$ws.Range("K132").Value = 8565.414000000001
$ws.Range("L132").Value = 14530.7139
$ws.Range("M132").Value = -6035.414000000001
$ws.Range("N132").Value = -19590.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4404.0386
$ws.Range("I7").Value = 3928.5715
$ws.Range("J7").Value = 6401
$ws.Range("K7").Value = 3928.5715
$ws.Range("L7").Value = 6401
$ws.Range("M7").Value = -3816.5715
$ws.Range("N7").Value = -6625
$ws.Range("H40").Value = 334869.34
$ws.Range("I40").Value = 334869.34
$ws.Range("K40").Value = 334869.34
$ws.Range("M40").Value = -334733.34
$ws.Range("H61").Value = 2828.8
$ws.Range("I61").Value = 947
$ws.Range("J61").Value = 4083.3333
$ws.Range("K61").Value = 947
$ws.Range("L61").Value = 4083.3333
$ws.Range("M61").Value = -745
$ws.Range("N61").Value = -4487.3333
$ws.Range("H100").Value = 2334.7693
$ws.Range("I100").Value = 2306.889
$ws.Range("J100").Value = 2397.5
$ws.Range("K100").Value = 2306.889
$ws.Range("L100").Value = 2397.5
$ws.Range("M100").Value = -1765.889
$ws.Range("N100").Value = -3479.5
$ws.Range("H113").Value = 2828.8
$ws.Range("I113").Value = 947
$ws.Range("J113").Value = 4083.3333
$ws.Range("K113").Value = 947
$ws.Range("L113").Value = 4083.3333
$ws.Range("M113").Value = 1223
$ws.Range("N113").Value = -8423.3333
$ws.Range("H122").Value = 3647.077
$ws.Range("I122").Value = 2235.3333
$ws.Range("J122").Value = 4857.143
$ws.Range("K122").Value = 6705.999899999999
$ws.Range("L122").Value = 14571.429
$ws.Range("M122").Value = -4255.999899999999
$ws.Range("N122").Value = -19471.429
$ws.Range("H126").Value = 4404.0386
$ws.Range("I126").Value = 3928.5715
$ws.Range("J126").Value = 6401
$ws.Range("K126").Value = 11785.7145
$ws.Range("L126").Value = 19203
$ws.Range("M126").Value = -9315.7145
$ws.Range("N126").Value = -24143
$ws.Range("H132").Value = 2977.6
$ws.Range("I132").Value = 2344.5
$ws.Range("J132").Value = 3927.25
$ws.Range("K132").Value = 7033.5
$ws.Range("L132").Value = 11781.75
$ws.Range("M132").Value = -4503.5
$ws.Range("N132").Value = -16841.75
$ws.Range("H136").Value = 1996.0555
$ws.Range("I136").Value = 2429.1428
$ws.Range("J136").Value = 1720.4546
$ws.Range("K136").Value = 7287.428400000001
$ws.Range("L136").Value = 5161.3638
$ws.Range("M136").Value = -4737.428400000001
$ws.Range("N136").Value = -10261.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4635
$ws.Range("I96").Value = 2695
$ws.Range("J96").Value = 6575
$ws.Range("K96").Value = 2695
$ws.Range("L96").Value = 6575
$ws.Range("M96").Value = -1322
$ws.Range("N96").Value = -9321
$ws.Range("H113").Value = 989.875
$ws.Range("I113").Value = 1021.63635
$ws.Range("J113").Value = 920
$ws.Range("K113").Value = 3064.90905
$ws.Range("L113").Value = 2760
$ws.Range("M113").Value = -894.9090500000002
$ws.Range("N113").Value = -7100
$ws.Range("H126").Value = 1000
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -530
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 6175801
$ws.Range("I132").Value = 4218.091
$ws.Range("K132").Value = 12654.273
$ws.Range("M132").Value = -10124.273
$ws.Range("H136").Value = 3405.1428
$ws.Range("I136").Value = 2777.95
$ws.Range("J136").Value = 4973.125
$ws.Range("K136").Value = 8333.849999999999
$ws.Range("L136").Value = 14919.375
$ws.Range("M136").Value = -5783.849999999999
$ws.Range("N136").Value = -20019.375

